# Updates the cryptos list (prices / 1h volume %) per the latest data refresh.
# Note: some "Price" values look like plain numbers (e.g. "554.96"); a leading
# apostrophe is used for those so Excel stores them as text (matching the
# original inline-string/text formatting) instead of auto-converting to a
# numeric value and dropping significant trailing digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.916.66'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '2.424.02'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''554.96'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").Value = '''138.15'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +4.77%  '
$ws.Range("D10").Value = '''5.80'
$ws.Range("E10").Value = '  +3.70%  '
$ws.Range("D11").Value = '''0.361'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("D13").Value = '''24.69'
$ws.Range("E13").Value = '  +3.56%  '
$ws.Range("D14").Value = '2.853.90'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").Value = '59.765.09'
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("E16").Value = '  +4.78%  '
$ws.Range("D17").Value = '2.405.02'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").Value = '''11.38'
$ws.Range("E18").Value = '  +6.04%  '
$ws.Range("E19").Value = '  +4.16%  '
$ws.Range("D20").Value = '''334.87'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").Value = '''6.94'
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '''64.72'
$ws.Range("E23").Value = '  +2.90%  '
$ws.Range("D24").Value = '''0.170'
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").Value = '''8.65'
$ws.Range("E25").Value = '  +2.01%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").Value = '0.0₃0793'
$ws.Range("E28").Value = '  +7.42%  '
$ws.Range("E29").Value = '  +2.97%  '
$ws.Range("D30").Value = '''170.39'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("D31").Value = '''6.29'
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("D32").Value = '''18.69'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").Value = '''1.02'
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D35").Value = '''1.32'
$ws.Range("E35").Value = '  +5.23%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = '''40.12'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  +11.42%  '
$ws.Range("D41").Value = '''312.60'
$ws.Range("E41").Value = '  +6.40%  '
$ws.Range("D42").Value = '''3.75'
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("D43").Value = '''142.56'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("E45").Value = '  +4.00%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''19.23'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("B47").Value = 'Polygon'
$ws.Range("C47").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D47").Value = '''0.409'
$ws.Range("E47").Value = '  +6.02%  '
$ws.Range("D48").Value = '''0.572'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("E49").Value = '  +2.97%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("E51").Value = '  +5.16%  '
